$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "Jun_13" / "Jun_10" columns
# (old column B -> D, old column C -> E), making room for two newer
# MarketBeat snapshot columns ("Jun_15" and "Jun_17").
$ws.Range("B1:C1").EntireColumn.Insert()

# New header row: most-recent date first.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new columns with the same placeholder ("UN") used throughout the
# rest of the table for rows with no special rating-change note.
$ws.Range("B2:C27").Value = "UN"

# Match the narrow 8-character width already used for the date columns.
$ws.Columns("C").ColumnWidth = 7.1667
$ws.Columns("D").ColumnWidth = 7.1667
$ws.Columns("E").ColumnWidth = 7.1667
